# Update Active_Outages.xlsx - 6/16/2025, 10:52:48 AM
# Target sheet is "R1" (sheet1.xml), which is the workbook's first / active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing values -------------------------------------------------

# G2: "3876.0" -> "3876.1" (value is stored as text, not a number, in the source
# file, so force text formatting for the write then restore the cell style so
# no stray formatting is left behind).
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "3876.1"
$ws.Range("G2").Style = "Normal"

# G3: "15.5" -> "15.7" (also text)
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "15.7"
$ws.Range("G3").Style = "Normal"

# D4: "asq0342" -> "JED0123"
$ws.Range("D4").Value = "JED0123"

# --- Append new outage row (row 5) ------------------------------------------
# Mirrors the shape of row 4: a new in-progress outage on hub site JED0123.

$ws.Range("A5").Formula = "="""""
$ws.Range("B5").Value = "R4"
$ws.Range("C5").Formula = "="""""
$ws.Range("D5").Value = "JED0123"
$ws.Range("E5").Formula = "="""""
$ws.Range("F5").Formula = "="""""
$ws.Range("G5").Formula = "="""""
$ws.Range("H5").Formula = "="""""
$ws.Range("I5").Value = "SCECO"
$ws.Range("J5").Value = "In progress"
$ws.Range("K5").Formula = "="""""
$ws.Range("L5").Value = "Latis"
